# Apply AnimalEcology edits:
#  - update sample rows 5-8 (generation interval / count / fade time / movement radius / speed / rest time / generation areas)
#  - delete the now-redundant duplicate rows 9:15 (rows shift up by 7)
#  - refresh the view selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 5 (TestAnimalName0001) ---
$ws.Range("D5").Value = "29|33|39"
$ws.Range("F5").Value = 300000
$ws.Range("G5").Value = 8
$ws.Range("H5").Value = 300000
$ws.Range("I5").Value = 5000
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = "5|3|6"

# --- Row 6 (TestAnimalName0002) ---
$ws.Range("D6").Value = "30|34|40"
$ws.Range("F6").Value = 300000
$ws.Range("G6").Value = 8
$ws.Range("H6").Value = 300000
$ws.Range("I6").Value = 5000
$ws.Range("J6").Value = "100|300|500"
$ws.Range("K6").Value = "3|4|5"

# --- Row 7 (TestAnimalName0003) ---
$ws.Range("D7").Value = "31|35|41"
$ws.Range("F7").Value = 300000
$ws.Range("G7").Value = 8
$ws.Range("H7").Value = 300000
$ws.Range("I7").Value = 5000
$ws.Range("J7").Value = "100|300|500"
$ws.Range("K7").Value = "3|2|4"

# --- Row 8 (TestAnimalName0004) ---
$ws.Range("F8").Value = 300000
$ws.Range("G8").Value = 8
$ws.Range("H8").Value = 300000
$ws.Range("I8").Value = 5000
$ws.Range("J8").Value = "100|300|500"
$ws.Range("K8").Value = "2|3|4"

# Remove the duplicate placeholder rows 9:15 entirely (everything below shifts up by 7)
$ws.Range("A9:A15").EntireRow.Delete()

# Refresh selection/view like the saved workbook
$ws.Range("M7").Select()
